$d = $word.ActiveDocument

$para = $d.Paragraphs(1).Range
$insertPoint = $d.Range($para.Start, $para.End - 1)
$insertPoint.Collapse(0)
$insertPoint.InsertAfter(" (")
$insertPoint.Collapse(0)
$insertPoint.InsertAfter("Changed main")
$insertPoint.Collapse(0)
$insertPoint.InsertAfter(")")
